# Update crypto price/volume figures per the Thu Apr 11 06:26:59 UTC 2024 GitHub Actions refresh.
# Rows 2-50 (1-based) hold one coin each; columns are A=rank, B=name, C=link, D=price, E=1h volume%.
# Numeric-looking price strings are entered with a leading apostrophe so Excel keeps them as literal
# text (matching the source data, which preserves things like "1.00" or thousand-dot formats).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.684.52"
$ws.Range("E2").Value = "  +1.88%  "

$ws.Range("D3").Value = "3.564.49"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'607.12"
$ws.Range("E5").Value = "  +4.19%  "

$ws.Range("D6").Value = "'173.86"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("D7").Value = "3.559.32"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").Value = "  +3.25%  "

$ws.Range("D11").Value = "'7.47"
$ws.Range("E11").Value = "  +9.89%  "

$ws.Range("D12").Value = "'0.588"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "'46.81"
$ws.Range("E13").Value = "  -1.72%  "

$ws.Range("E14").Value = "  +0.42%  "

$ws.Range("D15").Value = "4.149.15"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").Value = "'8.42"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "'615.70"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").Value = "3.573.24"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").Value = "70.801.68"
$ws.Range("E19").Value = "  +2.08%  "

$ws.Range("E20").Value = "  -2.23%  "

$ws.Range("D21").Value = "'17.42"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").Value = "'0.887"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").Value = "'9.43"
$ws.Range("E23").Value = "  -16.44%  "

$ws.Range("D24").Value = "'16.05"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "'97.40"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").Value = "'3.82"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").Value = "'33.51"
$ws.Range("E29").Value = "  +1.44%  "

$ws.Range("D30").Value = "'9.12"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("D31").Value = "'8.49"
$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("D32").Value = "'3.05"
$ws.Range("E32").Value = "  -3.82%  "

$ws.Range("D33").Value = "'7.02"
$ws.Range("E33").Value = "  -0.62%  "

$ws.Range("E34").Value = "  -2.68%  "

$ws.Range("D35").Value = "'610.18"
$ws.Range("E35").Value = "  -5.42%  "

$ws.Range("D36").Value = "'3.71"
$ws.Range("E36").Value = "  +5.31%  "

$ws.Range("E37").Value = "  -1.43%  "

$ws.Range("D38").Value = "'10.84"
$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("E39").Value = "  +5.49%  "

$ws.Range("D40").Value = "'57.31"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +3.25%  "

$ws.Range("D43").Value = "3.386.48"
$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("D44").Value = "'0.321"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'3.00"
$ws.Range("E45").Value = "  +8.05%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'33.12"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").Value = "0.0₃0708"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").Value = "'0.131"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").Value = "'132.84"
$ws.Range("E50").Value = "  +0.04%  "
